$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8038198253328799, 1.013021993066132, 0.8101406104434412, 2.5904784035739254)"
$ws.Range("C2").Value = "NIG(0.8812567081071281, 0.6102771038003127, 4.502810053177122, 5.0309814430203765)"
$ws.Range("D2").Value = "MIE(7.261820786819429, 4.3597126760473985, -7.277666547697973, 9.317270335800927)"
$ws.Range("E2").Value = "NIG(1.443837055104372, 1.189431416004987, 3.7117874615821096, 6.010673505955795)"
